# Append "-from develop brand" after the existing "-test modify" run, split
# around the pre-existing "_GoBack" bookmark so the bookmark keeps sitting
# between the two newly typed chunks, exactly like:
#
#   ...-test modify-from <bookmark/>develop brand
#
# which renders as "...-test modify-from develop brand".

$d = $word.ActiveDocument

$fromText  = "-from "
$brandText = "develop brand"

# Locate the existing "-test modify" run; Find.Execute collapses $rng to the
# matched text, so $rng.End is the position right after "...modify".
$rng = $d.Content
$found = $rng.Find.Execute("-test modify")
if (-not $found) {
    throw "Could not find '-test modify' in the document"
}
$afterTestModify = $rng.End

# Insert "-from " immediately after "-test modify". Using InsertBefore on a
# range collapsed at that position places the new text ahead of the
# "_GoBack" bookmark that already sits there (it merges into the preceding
# run, inheriting its Arial/23282B/kern36/sz27 formatting).
$gap1 = $d.Range($afterTestModify, $afterTestModify)
$gap1.InsertBefore($fromText)

# Position right after "-from " text, i.e. right after the bookmark.
$afterFrom = $afterTestModify + $fromText.Length

# Insert "develop brand" after that point (InsertAfter on a collapsed range
# here lands the new text following the bookmark, matching the diff order).
$gap2 = $d.Range($afterFrom, $afterFrom)
$gap2.InsertAfter($brandText)

# Explicitly format the newly created "develop brand" run to match the
# surrounding run's character formatting (Arial / Times New Roman east-asian
# fallback / color 23282B / kerning 18pt / size 13.5pt).
#
# NOTE: Font.NameBi is deliberately NOT touched here -- in this COM shim it
# has a destructive side effect that rewrites the "w:cs" font of unrelated
# runs elsewhere in the document, so it is avoided even though it would
# otherwise be the most direct way to stamp <w:rFonts w:cs="Arial"/>.
$newRun = $d.Range($afterFrom, $afterFrom + $brandText.Length)
$newRun.Font.Name        = "Arial"
$newRun.Font.NameFarEast = "Times New Roman"
$newRun.Font.NameOther   = "Arial"
$newRun.Font.Color       = 2828323
$newRun.Font.Kerning     = 18
$newRun.Font.Size        = 13.5

Write-Output $d.Paragraphs(1).Range.Text
